# Suppress any confirmation dialogs (e.g. sheet delete prompts)
$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# --- First sheet: rename existing Sheet1 -> TestValidLogin and fill in data ---
$ws1 = $wb.ActiveSheet
$ws1.Name = "TestValidLogin"

# Data typed first (rows of credentials), then headers added above them -
# this reproduces the shared-string ordering (admin, manager, Username, Password, ...)
$ws1.Range("A2").Value = "admin"
$ws1.Range("B2").Value = "manager"
$ws1.Range("A1").Value = "Username"
$ws1.Range("B1").Value = "Password"

# Auto-fit the two columns to their content
$ws1.Columns("A:A").AutoFit()
$ws1.Columns("B:B").AutoFit()

# Leave the cursor on B2 after finishing this sheet
[void]$ws1.Range("B2").Select()

# --- Second sheet: add a new sheet after the first one ---
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1)
$ws2.Name = "TestInvalidLogin"

# Headers typed first, then the invalid credential data below them
$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "passowrd"
$ws2.Range("A2").Value = "abcd"
$ws2.Range("B2").Value = "xyz"

# Leave the cursor one row below the data, on this (now active) sheet
[void]$ws2.Range("B3").Select()

# Make the second (invalid-login) sheet the active tab, matching the
# workbook's activeTab="1" / tabSelected behaviour
[void]$ws2.Activate()
